$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DropItems column (F) now stores item IDs instead of item names.
# Rows 2 and 5 hold a single numeric id; rows 3 and 4 hold a ";"-joined
# text list of ids.
$ws.Range("F2").Value = 1001
$ws.Range("F3").Value = "1001;1003"
$ws.Range("F4").Value = "1001;1003;1005"
$ws.Range("F5").Value = 1007

# Column width adjustments (values chosen so the pixel-quantised result the
# Excel object model actually persists lands as close as possible to the
# target widths of 11.25 / 33.25 / 23.75 characters)
$ws.Range("E1").ColumnWidth = 10.55
$ws.Range("F1").ColumnWidth = 32.55
$ws.Range("G1").ColumnWidth = 23

# Update the active selection
$ws.Range("F3").Select() | Out-Null
